$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new shared-string values in the exact order they were first
# introduced in the target workbook, so the generated sharedStrings.xml
# table lines up with the authoritative diff (index 12 .. 25).

# index 12
$ws.Range("D4").Value = "{'max_depth': 300, 'max_features': 6, 'n_estimators': 300}"
# index 13
$ws.Range("D5").Value = "{'max_depth': 200, 'max_features': 2, 'n_estimators': 100}"
# index 14
$ws.Range("D6").Value = "{'max_depth': 110, 'max_features': 2, 'n_estimators': 200}"
# index 15
$ws.Range("D7").Value = "{'max_depth': 300, 'max_features': 2, 'n_estimators': 200}"
# index 16
$ws.Range("D9").Value = "{'max_depth': 90, 'max_features': 6, 'n_estimators': 100}"
# index 17
$ws.Range("D10").Value = " {'max_depth': 80, 'max_features': 2, 'n_estimators': 300}"
# index 18
$ws.Range("D8").Value = "{'max_depth': 150, 'max_features': 4, 'n_estimators': 100}"
# index 19
$ws.Range("A12").Value = "Hypertension_clinical"
# index 20
$ws.Range("D12").Value = "{'max_depth': 200, 'max_features': 2, 'n_estimators': 300}"
# index 21
$ws.Range("D13").Value = "{'max_depth': 150, 'max_features': 2, 'n_estimators': 200}"
# index 22
$ws.Range("D14").Value = "{'max_depth': 110, 'max_features': 4, 'n_estimators': 200}"
# index 23
$ws.Range("D16").Value = "{'max_depth': 100, 'max_features': 3, 'n_estimators': 100}"
# index 24
$ws.Range("D18").Value = "{'max_depth': 200, 'max_features': 2, 'n_estimators': 200}"
# index 25
$ws.Range("D19").Value = "{'max_depth': 90, 'max_features': 2, 'n_estimators': 100}"

# Reused existing strings
$ws.Range("D15").Value = "{'max_depth': 110, 'max_features': 2, 'n_estimators': 200}"
$ws.Range("D17").Value = "{'max_depth': 90, 'max_features': 3, 'n_estimators': 100}"

# --- Numeric / remaining text fills, row by row ---

# Row 4 (B4=7 already set)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 91.379
$ws.Range("G4").Value = 0

# Row 5 (B5=8 already set)
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 93.103
$ws.Range("G5").Value = 0

# Row 6 (B6=9 already set)
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 86.207
$ws.Range("G6").Value = 3.448

# Row 7 (B7=10 already set)
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 92.24
$ws.Range("G7").Value = 2.096

# Row 8 (B8=11 already set)
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 93.927
$ws.Range("G8").Value = 3.318

# Row 9: B9 changes 12 -> 20
$ws.Range("B9").Value = 20
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 95.69
$ws.Range("G9").Value = 0.862

# Row 10: B10 changes 20 -> 50; C10 stays "Random Forest"
$ws.Range("B10").Value = 50
$ws.Range("C10").Value = "Random Forest"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 95.69
$ws.Range("G10").Value = 2.586

# Row 11: removed entirely (was A11=Hypertension_trans,B11=50,C11=Random Forest)
$ws.Range("A11:G11").ClearContents()

# New rows 12-19: Hypertension_clinical dataset
$ws.Range("B12").Value = 5
$ws.Range("C12").Value = "Random Forest"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 83.56
$ws.Range("G12").Value = 6.248

$ws.Range("A13").Value = "Hypertension_clinical"
$ws.Range("B13").Value = 6
$ws.Range("C13").Value = "Random Forest"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 82.087
$ws.Range("G13").Value = 7.292

$ws.Range("A14").Value = "Hypertension_clinical"
$ws.Range("B14").Value = 7
$ws.Range("C14").Value = "Random Forest"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 83.576
$ws.Range("G14").Value = 4.487

$ws.Range("A15").Value = "Hypertension_clinical"
$ws.Range("B15").Value = 8
$ws.Range("C15").Value = "Random Forest"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 83.927
$ws.Range("G15").Value = 5.73

$ws.Range("A16").Value = "Hypertension_clinical"
$ws.Range("B16").Value = 9
$ws.Range("C16").Value = "Random Forest"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 84.293
$ws.Range("G16").Value = 5.212

$ws.Range("A17").Value = "Hypertension_clinical"
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = "Random Forest"
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 85.579
$ws.Range("G17").Value = 5.537

$ws.Range("A18").Value = "Hypertension_clinical"
$ws.Range("B18").Value = 11
$ws.Range("C18").Value = "Random Forest"
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 85.051
$ws.Range("G18").Value = 5.034

$ws.Range("A19").Value = "Hypertension_clinical"
$ws.Range("B19").Value = 12
$ws.Range("C19").Value = "Random Forest"
$ws.Range("E19").Value = 3
$ws.Range("F19").Value = 84.289
$ws.Range("G19").Value = 3.679

# --- Column widths ---
$ws.Range("A1").EntireColumn.ColumnWidth = 19.85546875
$ws.Range("B1").EntireColumn.ColumnWidth = 12.7109375
$ws.Range("C1").EntireColumn.ColumnWidth = 14.7109375
$ws.Range("D1").EntireColumn.ColumnWidth = 23.5703125
$ws.Range("E1").EntireColumn.ColumnWidth = 6.140625
$ws.Range("F1").EntireColumn.ColumnWidth = 8.28515625
$ws.Range("G1").EntireColumn.ColumnWidth = 6.42578125

# --- Sheet view: scroll + selection ---
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("H13").Select()
